$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 50,4

$arr[0,0] = 'Bitcoin'
$arr[0,1] = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$arr[0,2] = '29.464.96'
$arr[0,3] = '  +1.91%  '
$arr[1,0] = 'Ethereum'
$arr[1,1] = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$arr[1,2] = '1.983.81'
$arr[1,3] = '  +5.58%  '
$arr[2,0] = 'TetherUSD'
$arr[2,1] = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$arr[2,2] = '''1.002'
$arr[2,3] = '  +0.07%  '
$arr[3,0] = 'BNB'
$arr[3,1] = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$arr[3,2] = '''327.61'
$arr[3,3] = '  +0.79%  '
$arr[4,0] = 'USDC'
$arr[4,1] = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$arr[4,2] = '''1.000'
$arr[4,3] = '  -0.05%  '
$arr[5,0] = 'XRP'
$arr[5,1] = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$arr[5,2] = '''0.4673'
$arr[5,3] = '  +1.69%  '
$arr[6,0] = 'Cardano'
$arr[6,1] = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$arr[6,2] = '''0.3936'
$arr[6,3] = '  +1.52%  '
$arr[7,0] = 'Dogecoin'
$arr[7,1] = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$arr[7,2] = '''0.07948'
$arr[7,3] = '  +0.96%  '
$arr[8,0] = 'Polygon'
$arr[8,1] = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$arr[8,2] = '''1.000'
$arr[8,3] = '  +1.47%  '
$arr[9,0] = 'Solana'
$arr[9,1] = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$arr[9,2] = '''22.83'
$arr[9,3] = '  +4.78%  '
$arr[10,0] = 'WrappedEther'
$arr[10,1] = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$arr[10,2] = '1.955.84'
$arr[10,3] = '  +4.29%  '
$arr[11,0] = 'Chainlink'
$arr[11,1] = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$arr[11,2] = '''7.246'
$arr[11,3] = '  +3.79%  '
$arr[12,0] = 'Polkadot'
$arr[12,1] = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$arr[12,2] = '''5.865'
$arr[12,3] = '  +3.84%  '
$arr[13,0] = 'TRON'
$arr[13,1] = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$arr[13,2] = '''0.07129'
$arr[13,3] = '  +2.40%  '
$arr[14,0] = 'Litecoin'
$arr[14,1] = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$arr[14,2] = '''88.61'
$arr[14,3] = '  +0.67%  '
$arr[15,0] = 'BinanceUSD'
$arr[15,1] = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$arr[15,2] = '''1.003'
$arr[15,3] = '  +0.18%  '
$arr[16,0] = 'ShibaInu'
$arr[16,1] = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$arr[16,2] = '''0.000009946'
$arr[16,3] = '  -0.25%  '
$arr[17,0] = 'Avalanche'
$arr[17,1] = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$arr[17,2] = '''17.32'
$arr[17,3] = '  +2.04%  '
$arr[18,0] = 'Dai'
$arr[18,1] = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$arr[18,2] = '''1.002'
$arr[18,3] = '  +0.14%  '
$arr[19,0] = 'WrappedBTC'
$arr[19,1] = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$arr[19,2] = '29.589.69'
$arr[19,3] = '  +2.32%  '
$arr[20,0] = 'Uniswap'
$arr[20,1] = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$arr[20,2] = '''5.533'
$arr[20,3] = '  +5.52%  '
$arr[21,0] = 'Cosmos'
$arr[21,1] = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$arr[21,2] = '''11.25'
$arr[21,3] = '  +2.64%  '
$arr[22,0] = 'WrappedliquidstakedEther2.0'
$arr[22,1] = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$arr[22,2] = '2.236.76'
$arr[22,3] = '  +6.28%  '
$arr[23,0] = 'Toncoin'
$arr[23,1] = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$arr[23,2] = '''2.112'
$arr[23,3] = '  +0.54%  '
$arr[24,0] = 'Monero'
$arr[24,1] = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$arr[24,2] = '''157.46'
$arr[24,3] = '  +0.78%  '
$arr[25,0] = 'EthereumClassic'
$arr[25,1] = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$arr[25,2] = '''19.64'
$arr[25,3] = '  +1.74%  '
$arr[26,0] = 'InternetComputer(DFINITY)'
$arr[26,1] = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$arr[26,2] = '''6.000'
$arr[26,3] = '  -0.36%  '
$arr[27,0] = 'BitcoinCash'
$arr[27,1] = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$arr[27,2] = '''120.14'
$arr[27,3] = '  +2.40%  '
$arr[28,0] = 'LidoDAOToken'
$arr[28,1] = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$arr[28,2] = '''1.960'
$arr[28,3] = '  +1.56%  '
$arr[29,0] = 'Stellar'
$arr[29,1] = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$arr[29,2] = '''0.09437'
$arr[29,3] = '  +1.25%  '
$arr[30,0] = 'ImmutableX'
$arr[30,1] = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$arr[30,2] = '''0.8903'
$arr[30,3] = '  -1.34%  '
$arr[31,0] = 'Filecoin'
$arr[31,1] = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$arr[31,2] = '''5.271'
$arr[31,3] = '  +0.33%  '
$arr[32,0] = 'ARBITRUM'
$arr[32,1] = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$arr[32,2] = '''1.343'
$arr[32,3] = '  +1.98%  '
$arr[33,0] = 'HuobiToken'
$arr[33,1] = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$arr[33,2] = '''3.181'
$arr[33,3] = '  -2.22%  '
$arr[34,0] = 'Hedera'
$arr[34,1] = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$arr[34,2] = '''0.05835'
$arr[34,3] = '  +1.34%  '
$arr[35,0] = 'TrustWalletToken'
$arr[35,1] = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$arr[35,2] = '''1.174'
$arr[35,3] = '  -1.08%  '
$arr[36,0] = 'VeChain'
$arr[36,1] = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$arr[36,2] = '''0.02124'
$arr[36,3] = '  +2.62%  '
$arr[37,0] = 'PEPE'
$arr[37,1] = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$arr[37,2] = '''0.000003330'
$arr[37,3] = '  +108.06%  '
$arr[38,0] = 'FraxShare'
$arr[38,1] = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$arr[38,2] = '''7.901'
$arr[38,3] = '  +3.07%  '
$arr[39,0] = 'TheSandbox'
$arr[39,1] = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$arr[39,2] = '''0.5751'
$arr[39,3] = '  +1.85%  '
$arr[40,0] = 'Algorand'
$arr[40,1] = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$arr[40,2] = '''0.1824'
$arr[40,3] = '  +3.40%  '
$arr[41,0] = 'Aptos'
$arr[41,1] = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$arr[41,2] = '''9.806'
$arr[41,3] = '  +1.49%  '
$arr[42,0] = 'EnergySwap'
$arr[42,1] = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$arr[42,2] = '''12.13'
$arr[42,3] = '  +2.00%  '
$arr[43,0] = 'Decentraland'
$arr[43,1] = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$arr[43,2] = '''0.5370'
$arr[43,3] = '  +0.48%  '
$arr[44,0] = 'MXToken'
$arr[44,1] = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$arr[44,2] = '''2.674'
$arr[44,3] = '  +6.38%  '
$arr[45,0] = 'RenderToken'
$arr[45,1] = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$arr[45,2] = '''2.164'
$arr[45,3] = '  -4.28%  '
$arr[46,0] = 'Cronos'
$arr[46,1] = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$arr[46,2] = '''0.06955'
$arr[46,3] = '  -1.24%  '
$arr[47,0] = 'NEARProtocol'
$arr[47,1] = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$arr[47,2] = '''1.865'
$arr[47,3] = '  +1.02%  '
$arr[48,0] = 'Quant'
$arr[48,1] = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$arr[48,2] = '''114.27'
$arr[48,3] = '  +1.13%  '
$arr[49,0] = 'WOONetwork'
$arr[49,1] = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$arr[49,2] = '''0.3096'
$arr[49,3] = '  +8.55%  '

$ws.Range("B2:E51").Value = $arr
